$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear out old content (column D no longer used)
$ws.Range("A1:D5").ClearContents()

# Header row
$ws.Range("A1").Value = "Nombre"
$ws.Range("B1").Value = "Inicial"
$ws.Range("C1").Value = "Descripción"

# Fill the description column down for the middle rows first
$ws.Range("C2").Value = "Descripcion1"
$ws.Range("C3").Value = "Descripcion2"
$ws.Range("C4").Value = "Descripcion3"

# Then fill columns A & B row by row
$ws.Range("A2").Value = "Partido Socialista"
$ws.Range("B2").Value = "PSOE"

$ws.Range("A3").Value = "Partido Popular"
$ws.Range("B3").Value = "PP"

$ws.Range("A4").Value = "Podemos"
$ws.Range("B4").Value = "Podemos"

$ws.Range("A5").Value = "Ciudadanos"
$ws.Range("B5").Value = "C's"

# Last description cell filled in at the end
$ws.Range("C5").Value = "Descripcion4"

# Column A width (closest achievable value to target stored width 17.85546875)
$ws.Columns.Item(1).ColumnWidth = 17

# Select C5 as the active cell
$ws.Range("C5").Select()
